$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("R2").Value = 1.33

# Row 3
$ws.Range("M3").Value = 1.05
$ws.Range("N3").Value = 8
$ws.Range("O3").Value = 1.41
$ws.Range("P3").Value = 2.62
$ws.Range("Q3").Value = 2.35
$ws.Range("R3").Value = 1.57

# Row 4
$ws.Range("M4").Value = 1.07
$ws.Range("O4").Value = 1.41
$ws.Range("P4").Value = 2.62
